# Build v2.1.2: Fix SearchCriteria variants and Schemas sheet grouping/sorting
#
# Summary of changes (per the OOXML diff):
#  - "Body"  sheet: row 3 becomes the revokeLcr.230301Request schema ref; rows 4-7 removed.
#  - "200"   sheet: row 3 becomes the revokeLcr.230301Response schema ref; rows 4-5 removed.
#  - "204"   sheet: new row 3 added for the revokeLcr.230301Response schema ref.
#  - "400"   sheet: row 3 becomes the errorResponse schema ref; rows 4-6 removed.
#  - "401","403","404","429","500" sheets: new row 3 added for the errorResponse1 schema ref.

$wb = $excel.ActiveWorkbook

# NOTE: this interpreter's function dispatch only reliably supports
# *positional* parameters - named parameters (`-ws ... -name ...`) bind
# incorrectly, so the helper below takes plain positional args.
function Set-SchemaRow3($ws, $section, $name) {
    $ws.Range("A3").Value = $section
    $ws.Range("B3").Value = $name
    $ws.Range("C3").Value = ""
    $ws.Range("D3").Value = ""
    $ws.Range("E3").Value = "schema"
    $ws.Range("F3").Value = ""
    $ws.Range("G3").Value = $name
    $ws.Range("H3").Value = ""
    $ws.Range("I3").Value = "Yes"
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = ""
    $ws.Range("L3").Value = ""
    $ws.Range("M3").Value = ""
    $ws.Range("N3").Value = ""
    $ws.Range("O3").Value = ""
}

# ---- "Body" sheet: replace row 3 with the revokeLcr.230301Request schema ref, drop rows 4-7 ----
$wsBody = $wb.Worksheets.Item("Body")
Set-SchemaRow3 $wsBody "body" "revokeLcr.230301Request"
$wsBody.Range("A4:O7").Clear()

# ---- "200" sheet: replace row 3 with the revokeLcr.230301Response schema ref, drop rows 4-5 ----
$ws200 = $wb.Worksheets.Item("200")
Set-SchemaRow3 $ws200 "content" "revokeLcr.230301Response"
$ws200.Range("A4:O5").Clear()

# ---- "204" sheet: add row 3 for the revokeLcr.230301Response schema ref ----
$ws204 = $wb.Worksheets.Item("204")
Set-SchemaRow3 $ws204 "content" "revokeLcr.230301Response"

# ---- "400" sheet: replace row 3 with the errorResponse schema ref, drop rows 4-6 ----
$ws400 = $wb.Worksheets.Item("400")
Set-SchemaRow3 $ws400 "content" "errorResponse"
$ws400.Range("A4:O6").Clear()

# ---- "401" sheet: add row 3 for the errorResponse1 schema ref ----
$ws401 = $wb.Worksheets.Item("401")
Set-SchemaRow3 $ws401 "content" "errorResponse1"

# ---- "403" sheet: add row 3 for the errorResponse1 schema ref ----
$ws403 = $wb.Worksheets.Item("403")
Set-SchemaRow3 $ws403 "content" "errorResponse1"

# ---- "404" sheet: add row 3 for the errorResponse1 schema ref ----
$ws404 = $wb.Worksheets.Item("404")
Set-SchemaRow3 $ws404 "content" "errorResponse1"

# ---- "429" sheet: add row 3 for the errorResponse1 schema ref ----
$ws429 = $wb.Worksheets.Item("429")
Set-SchemaRow3 $ws429 "content" "errorResponse1"

# ---- "500" sheet: add row 3 for the errorResponse1 schema ref ----
$ws500 = $wb.Worksheets.Item("500")
Set-SchemaRow3 $ws500 "content" "errorResponse1"

Write-Output "revokeLcr.230301 sheets updated"
